$d = $word.ActiveDocument
$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# Change 1: Collapse the three long "CORE COMPETENCIES" bullet paragraphs
# into a single summary paragraph.
# ---------------------------------------------------------------------------
$coreHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "CORE COMPETENCIES") {
        $coreHeading = $i
        break
    }
}

if ($coreHeading -ne $null) {
    $firstIdx = $coreHeading + 1
    # Replace the first of the three paragraphs with the condensed text.
    $firstPara = $d.Paragraphs.Item($firstIdx)
    $firstPara.Range.Text = "Product Management & Strategy " + $bullet + " Technical Product Development " + $bullet + " Platform & Infrastructure"

    # Remove the next two (now-redundant) detailed paragraphs.
    $d.Paragraphs.Item($firstIdx + 1).Range.Delete()
    $d.Paragraphs.Item($firstIdx + 1).Range.Delete()
}

# ---------------------------------------------------------------------------
# Change 2: Append a new "TECHNICAL SKILLS" section at the end of the body.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading 2"

$lines = @(
    "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics",
    "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration",
    "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training"
)

foreach ($line in $lines) {
    $prevRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
    $prevRange.Collapse(0)
    $prevRange.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Style = "Normal"
    $newPara.Range.Text = $line
}
